$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill")
$ws.Rows.Item(20).Insert()
Write-Host "Row 20 after insert:" $ws.Cells.Item(20,1).Value2
Write-Host "Row 21 after insert:" $ws.Cells.Item(21,1).Value2
